$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 5, shifting existing rows 5-7 down to 6-8
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5 with the new data
$ws.Range("A5").Value = 0.5
$ws.Range("B5").Value = 0.101042934281782
$ws.Range("C5").Value = 0.303128802845346
